$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# The data rows 2-5 (columns A, Q, R, AC) are cyclically shifted up by one row:
# new row2 <- old row3, new row3 <- old row4, new row4 <- old row5, new row5 <- old row2
# Capture the original values first so the shifting assignment does not clobber source data.

$origA = @{}
$origQ = @{}
$origR = @{}
$origAC = @{}

foreach ($r in 2..5) {
    $origA[$r]  = $ws.Cells.Item($r, 1).Value2    # column A
    $origQ[$r]  = $ws.Cells.Item($r, 17).Value2   # column Q
    $origR[$r]  = $ws.Cells.Item($r, 18).Value2   # column R
    $origAC[$r] = $ws.Cells.Item($r, 29).Value2   # column AC
}

$mapping = @{ 2 = 3; 3 = 4; 4 = 5; 5 = 2 }

foreach ($r in 2..5) {
    $src = $mapping[$r]
    $ws.Cells.Item($r, 1).Value2  = $origA[$src]
    $ws.Cells.Item($r, 17).Value2 = $origQ[$src]
    $ws.Cells.Item($r, 18).Value2 = $origR[$src]
    $ws.Cells.Item($r, 29).Value2 = $origAC[$src]
}

$wb.Save()
